$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-12 Thursday", 2) | Out-Null
$d.Content.Find.Execute("182÷5=36, 2", $true, $false, $false, $false, $false, $true, 1, $false, "363÷9=40, 3", 2) | Out-Null
$d.Content.Find.Execute("434÷6=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "568÷3=189, 1", 2) | Out-Null
$d.Content.Find.Execute("193÷9=21, 4", $true, $false, $false, $false, $false, $true, 1, $false, "763÷7=109, 0", 2) | Out-Null
$d.Content.Find.Execute("338÷5=67, 3", $true, $false, $false, $false, $false, $true, 1, $false, "150÷4=37, 2", 2) | Out-Null
$d.Content.Find.Execute("741÷8=92, 5", $true, $false, $false, $false, $false, $true, 1, $false, "953÷3=317, 2", 2) | Out-Null
$d.Content.Find.Execute("825÷6=137, 3", $true, $false, $false, $false, $false, $true, 1, $false, "932÷8=116, 4", 2) | Out-Null
$d.Content.Find.Execute("228÷6=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "845÷4=211, 1", 2) | Out-Null
$d.Content.Find.Execute("523÷8=65, 3", $true, $false, $false, $false, $false, $true, 1, $false, "399÷5=79, 4", 2) | Out-Null
$d.Content.Find.Execute("898÷4=224, 2", $true, $false, $false, $false, $false, $true, 1, $false, "378÷6=63, 0", 2) | Out-Null
$d.Content.Find.Execute("390÷3=130, 0", $true, $false, $false, $false, $false, $true, 1, $false, "692÷8=86, 4", 2) | Out-Null
$d.Content.Find.Execute("526÷2=263, 0", $true, $false, $false, $false, $false, $true, 1, $false, "972÷4=243, 0", 2) | Out-Null
$d.Content.Find.Execute("230÷6=38, 2", $true, $false, $false, $false, $false, $true, 1, $false, "119÷8=14, 7", 2) | Out-Null
$d.Content.Find.Execute("218÷3=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "847÷7=121, 0", 2) | Out-Null
$d.Content.Find.Execute("609÷8=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "838÷8=104, 6", 2) | Out-Null
$d.Content.Find.Execute("830÷6=138, 2", $true, $false, $false, $false, $false, $true, 1, $false, "937÷7=133, 6", 2) | Out-Null
$d.Content.Find.Execute("839÷4=209, 3", $true, $false, $false, $false, $false, $true, 1, $false, "687÷2=343, 1", 2) | Out-Null
$d.Content.Find.Execute("549÷8=68, 5", $true, $false, $false, $false, $false, $true, 1, $false, "926÷8=115, 6", 2) | Out-Null
$d.Content.Find.Execute("171÷3=57, 0", $true, $false, $false, $false, $false, $true, 1, $false, "160÷3=53, 1", 2) | Out-Null
$d.Content.Find.Execute("624÷6=104, 0", $true, $false, $false, $false, $false, $true, 1, $false, "470÷3=156, 2", 2) | Out-Null
$d.Content.Find.Execute("593÷6=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "212÷5=42, 2", 2) | Out-Null
$d.Content.Find.Execute("764÷2=382, 0", $true, $false, $false, $false, $false, $true, 1, $false, "585÷2=292, 1", 2) | Out-Null
$d.Content.Find.Execute("816÷2=408, 0", $true, $false, $false, $false, $false, $true, 1, $false, "805÷3=268, 1", 2) | Out-Null
$d.Content.Find.Execute("506÷7=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "668÷3=222, 2", 2) | Out-Null
$d.Content.Find.Execute("827÷5=165, 2", $true, $false, $false, $false, $false, $true, 1, $false, "586÷4=146, 2", 2) | Out-Null
$d.Content.Find.Execute("102÷7=14, 4", $true, $false, $false, $false, $false, $true, 1, $false, "258÷4=64, 2", 2) | Out-Null
